$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2 = @(3.182878228561681, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 6.048734245549538)
    3 = @(3.182878228561681, 1.65323645889881, 0.1529057820181812, 0.4998867070740569, 0, 5.488907176552729)
    4 = @(0.02258322285507441, 0.0001537489499301437, 3.082599426703578, 6.48142807727062, 0, 9.586764475779203)
    5 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 0.4998867070740569, 0, 8.418600821238126)
    6 = @(0.7287194209349384, 1.65323645889881, 0.7127328510149897, 6.48142807727062, 1, 9.576116808119359)
    7 = @(1.505614041169197, 1.65323645889881, 0.7127328510149897, 0.4998867070740569, 0, 4.371470058157054)
    8 = @(3.182878228561681, 1.65323645889881, 3.082599426703578, 6.48142807727062, 0, 14.40014219143469)
}

foreach ($r in $data.Keys) {
    $vals = $data[$r]
    $ws.Range("B$r").Value = $vals[0]
    $ws.Range("C$r").Value = $vals[1]
    $ws.Range("D$r").Value = $vals[2]
    $ws.Range("E$r").Value = $vals[3]
    $ws.Range("F$r").Value = $vals[4]
    $ws.Range("G$r").Value = $vals[5]
}
